# Apply updated crypto price/volume figures (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.557.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.82%  "

$ws.Range("D3").Value = "'3.161.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.45%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").Value = "'538.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.79%  "

$ws.Range("D6").Value = "'140.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.28%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("E8").Value = "  +10.88%  "

$ws.Range("D9").Value = "'7.34"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.63%  "

$ws.Range("D10").Value = "'0.110"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.30%  "

$ws.Range("D11").Value = "'0.423"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.75%  "

$ws.Range("E12").Value = "  +2.93%  "

$ws.Range("D13").Value = "'3.691.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.84%  "

$ws.Range("D14").Value = "'26.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.63%  "

$ws.Range("D15").Value = "'0.0000170"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.96%  "

$ws.Range("D16").Value = "'58.524.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.66%  "

$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").Value = "'6.25"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +7.12%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "'3.150.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.90%  "

$ws.Range("D19").Value = "'13.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.55%  "

$ws.Range("D20").Value = "'8.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.91%  "

$ws.Range("D21").Value = "'378.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +9.08%  "

$ws.Range("E22").Value = "  +0.13%  "

$ws.Range("D23").Value = "'5.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("D24").Value = "'70.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.38%  "

$ws.Range("D25").Value = "'0.519"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.78%  "

$ws.Range("E26").Value = "  +2.74%  "

$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.06%  "

$ws.Range("D28").Value = "'7.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +12.55%  "

$ws.Range("D29").Value = "'0.0₃0884"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.85%  "

$ws.Range("E30").Value = "  +3.25%  "

$ws.Range("D31").Value = "'6.19"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.46%  "

$ws.Range("D32").Value = "'21.92"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.51%  "

$ws.Range("D33").Value = "'5.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.85%  "

$ws.Range("E34").Value = "  +5.95%  "

$ws.Range("D35").Value = "'161.28"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.75%  "

$ws.Range("D36").Value = "'6.27"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.54%  "

$ws.Range("D37").Value = "'1.36"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +11.09%  "

$ws.Range("D38").Value = "'25.67"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.17%  "

$ws.Range("D39").Value = "'1.69"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.61%  "

$ws.Range("D40").Value = "'2.651.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.82%  "

$ws.Range("D41").Value = "'0.0683"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.95%  "

$ws.Range("D42").Value = "'4.24"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.07%  "

$ws.Range("D43").Value = "'38.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.72%  "

$ws.Range("D44").Value = "'0.705"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.05%  "

$ws.Range("D45").Value = "'0.0278"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.03%  "

$ws.Range("D46").Value = "'0.998"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.21%  "

$ws.Range("E47").Value = "  +13.25%  "

$ws.Range("D48").Value = "'6.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.12%  "

$ws.Range("D49").Value = "'0.982"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.16%  "

$ws.Range("D50").Value = "'20.27"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.01%  "

$ws.Range("D51").Value = "'0.754"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.02%  "
